# Apply "new results" edit to the my_sheet_si worksheet:
#  - drop the stray empty B3 cell
#  - append 10 new data rows (4-13) with the new result set
#  - sheet dimension grows from A1:H3 to A1:H13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: remove the stray empty B3 cell ---
$ws.Range("B3").ClearContents()

# --- New rows 4-13 of results ---

# Row 4
$ws.Range("A4").Value = 8
$ws.Range("C4").Value = 0.2872303247797749
$ws.Range("D4").Value = '0:00:15.920702'
$ws.Range("E4").Value = 0.2028611111111111
$ws.Range("F4").Value = 0.3715995384484387
$ws.Range("G4").Value = '[''MAPK12'', ''TRIM23'', ''CIT'', ''TNS1'', ''PIK3C2A'', ''RAF1'', ''PPP3R1'', ''MAPK1'']'

# Row 5
$ws.Range("A5").Value = 25
$ws.Range("C5").Value = 0.1876826590545108
$ws.Range("D5").Value = '0:00:16.956974'
$ws.Range("E5").Value = 0.16532
$ws.Range("F5").Value = 0.2100453181090215
$ws.Range("G5").Value = '[''AURKA'', ''NFKBIB'', ''PI3'', ''PIK3CB'', ''AKT1'', ''KIF22'', ''EXOC2'', ''PRKACB'', ''PIK3C3'', ''PLCG1'', ''DLGAP5'', ''RHOBTB1'', ''NFKBIA'', ''NFKB1'', ''GABARAPL2'', ''TLR3'', ''RAC3'', ''PPP3CA'', ''SOS2'', ''TNF'', ''CHEK1'', ''SYK'', ''RPS6KB1'', ''PPP1CB'', ''PRKCI'']'

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("C6").Value = 0.3888482498927541
$ws.Range("D6").Value = '0:00:15.162950'
$ws.Range("E6").Value = 0.4126000000000001
$ws.Range("F6").Value = 0.3650964997855081
$ws.Range("G6").Value = '[''POLE2'', ''POLR2H'', ''POLR2I'', ''NME3'', ''GTF2H4'']'

# Row 7
$ws.Range("A7").Value = 4
$ws.Range("C7").Value = 0.2473947613832676
$ws.Range("D7").Value = '0:00:15.366835'
$ws.Range("E7").Value = 0.185
$ws.Range("F7").Value = 0.3097895227665352
$ws.Range("G7").Value = '[''RPS6KA2'', ''MKKS'', ''MAPK12'', ''CDK16'']'

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("C8").Value = 0.417886679229579
$ws.Range("D8").Value = '0:00:15.682973'
$ws.Range("E8").Value = 0.5034285714285714
$ws.Range("F8").Value = 0.3323447870305865
$ws.Range("G8").Value = '[''PSMC6'', ''PSME2'', ''UBA52'', ''PSMA6'', ''PPP2CB'', ''NFKBIA'', ''PPP2R5B'']'

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("C9").Value = 0.3839870264497305
$ws.Range("D9").Value = '0:00:15.047124'
$ws.Range("E9").Value = 0.3473928571428572
$ws.Range("F9").Value = 0.4205811957566039
$ws.Range("G9").Value = '[''YWHAZ'', ''INS'', ''GAPDH'', ''CDKN1B'', ''AURKA'', ''PIN1'', ''NFKBIA'']'

# Row 10
$ws.Range("A10").Value = 20
$ws.Range("C10").Value = 0.2117678632919119
$ws.Range("D10").Value = '0:00:16.964040'
$ws.Range("E10").Value = 0.1725666666666666
$ws.Range("F10").Value = 0.2509690599171571
$ws.Range("G10").Value = '[''RALA'', ''MAPK12'', ''SOS2'', ''HGF'', ''NFKB1'', ''PPP1CA'', ''FGFR3'', ''RAC1'', ''TNS1'', ''KIF22'', ''PPP1CC'', ''PPP2R1A'', ''AURKA'', ''STAG2'', ''SOS1'', ''PPP2R5B'', ''DCN'', ''DLGAP5'', ''SIRT1'', ''PPP3CA'']'

# Row 11
$ws.Range("A11").Value = 5
$ws.Range("C11").Value = 0.3978440954282857
$ws.Range("D11").Value = '0:00:16.238540'
$ws.Range("E11").Value = 0.45
$ws.Range("F11").Value = 0.3456881908565715
$ws.Range("G11").Value = '[''ANAPC7'', ''UBE2V2'', ''SKP1'', ''NFKB1'', ''FBXL15'']'

# Row 12
$ws.Range("A12").Value = 8
$ws.Range("C12").Value = 0.4865680695869101
$ws.Range("D12").Value = '0:00:16.228390'
$ws.Range("E12").Value = 0.6094444444444443
$ws.Range("F12").Value = 0.3636916947293759
$ws.Range("G12").Value = '[''KALRN'', ''PIK3R2'', ''GNA11'', ''GNRHR'', ''GAST'', ''KISS1R'', ''GNG4'', ''GNG7'']'

# Row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = 0.05330980403894163
$ws.Range("D13").Value = '0:00:02.387701'
$ws.Range("E13").Value = 0.04619696969696969
$ws.Range("F13").Value = 0.06042263838091356
$ws.Range("G13").Value = '[''FCGR3B'', ''CLPS'', ''LILRB3'', ''TP53'', ''PDE5A'', ''IL3RA'', ''CD79A'', ''CHP2'', ''PIK3CD'', ''TNF'', ''C1DP3'']'

Write-Output ("Dimension now: " + $ws.UsedRange.Address())
